$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 117845141.8069585
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 117845790.261642

# Row 3
$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 29.84159230404497

# Row 4
$ws.Range("B4").Value = 0.003994804209775715
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 9844.520545567508
$ws.Range("E4").Value = 645.3272768299601
$ws.Range("G4").Value = 10490.16460759753
